$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split the "URL" column (C) into two columns: "BaseURL" and "pathQuery".
# Inserting before column C shifts the old C..N data to D..O, carrying
# over styles/number formats from the original column C.
$ws.Columns("C").Insert()

# Header row
$ws.Range("C1").Value = "BaseURL"
$ws.Range("D1").Value = "pathQuery"

# Data row: split the old full URL into base URL + path query
$ws.Range("C2").Value = "https://devents.azure-api.net/Event"
$ws.Range("D2").Value = "/GetParticipants"

# The new BaseURL cell (C2) should carry the same (hyperlink) style as the
# pathQuery cell next to it (which inherited the original URL cell's style).
$ws.Range("D2").Copy()
$ws.Range("C2").PasteSpecial(-4122)

# Restore the column width on the new BaseURL column to match the original
# URL column's width (now carried by pathQuery, column D).
$ws.Columns("C").ColumnWidth = $ws.Columns("D").ColumnWidth

# Match the view state captured in the saved workbook (scrolled right with
# N1 selected).
$ws.Activate()
$ws.Range("N1").Select()
